{"js": "// Update the single-column benchmark-stats table in place.\n// Rows are addressed 0-based, top to bottom, matching the table's row order.\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Simple single-value replacements.\nconst updates = [\n  [0, \"0M\"],\n  [1, \"0M\"],\n  [2, \"0M\"],\n  [3, \"1839\"],\n  [6, \"0.03652\"],\n  [7, \"0.01280\"],\n  [11, \"2.81095\"],\n  // These three rows previously packed a tab-separated 9-field stat dump into\n  // one cell; they now hold the single summary value that rows 0/1/2 used to\n  // show before those became \"0M\" placeholders.\n  [43, \"99.8\"],\n  [44, \"2.81\"],\n  [45, \"1401\"],\n];\n\nfor (const [rowIndex, newText] of updates) {\n  const cell = table.getCell(rowIndex, 0);\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Update the single-column benchmark-stats table in place.\n# Word COM tables/cells are 1-indexed, so row N (1-based) == 0-based row (N-1).\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"1839\"\n    7  = \"0.03652\"\n    8  = \"0.01280\"\n    12 = \"2.81095\"\n    # These three rows previously packed a tab-separated 9-field stat dump\n    # into one cell; they now hold the single summary value that rows 1/2/3\n    # used to show before those became \"0M\" placeholders.\n    44 = \"99.8\"\n    45 = \"2.81\"\n    46 = \"1401\"\n}\n\nforeach ($row in $updates.Keys) {\n    $t.Cell($row, 1).Range.Text = $updates[$row]\n}\n"}
